# Auto-generated edit script: updates Leve profit sheet values per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 2442.889
$ws.Range("I20").Value = 1566.1428
$ws.Range("K20").Value = 1566.1428
$ws.Range("M20").Value = -1336.1428

$ws.Range("H35").Value = 2442.889
$ws.Range("I35").Value = 1566.1428
$ws.Range("K35").Value = 1566.1428
$ws.Range("M35").Value = -1187.1428

$ws.Range("H70").Value = 169215
$ws.Range("J70").Value = 252325
$ws.Range("L70").Value = 756975
$ws.Range("N70").Value = -757515

$ws.Range("H73").Value = 169215
$ws.Range("J73").Value = 252325
$ws.Range("L73").Value = 756975
$ws.Range("N73").Value = -758847

$ws.Range("H76").Value = 77007900
$ws.Range("I76").Value = 91893.664
$ws.Range("K76").Value = 91893.664
$ws.Range("M76").Value = -91578.664

$ws.Range("H79").Value = 77007900
$ws.Range("I79").Value = 91893.664
$ws.Range("K79").Value = 91893.664
$ws.Range("M79").Value = -90801.664

$ws.Range("H80").Value = 5489.2383
$ws.Range("I80").Value = 10889
$ws.Range("J80").Value = 580.36365
$ws.Range("K80").Value = 32667
$ws.Range("L80").Value = 1741.09095
$ws.Range("M80").Value = -31669
$ws.Range("N80").Value = -3737.09095

$ws.Range("H83").Value = 5489.2383
$ws.Range("I83").Value = 10889
$ws.Range("J83").Value = 580.36365
$ws.Range("K83").Value = 98001
$ws.Range("L83").Value = 5223.27285
$ws.Range("M83").Value = -93009
$ws.Range("N83").Value = -15207.27285

$ws.Range("H113").Value = 6248.5
$ws.Range("I113").Value = 5831.3335
$ws.Range("K113").Value = 5831.3335
$ws.Range("M113").Value = -2577.3335

$ws.Range("H138").Value = 5338.4775
$ws.Range("I138").Value = 3279.4092
$ws.Range("J138").Value = 6004.647
$ws.Range("K138").Value = 9838.2276
$ws.Range("L138").Value = 18013.941
$ws.Range("M138").Value = -4698.2276
$ws.Range("N138").Value = -28293.941

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3691.4407
$ws.Range("I32").Value = 2649.9807
$ws.Range("K32").Value = 2649.9807
$ws.Range("M32").Value = -2362.9807

$ws.Range("H63").Value = 6287.75
$ws.Range("I63").Value = 5043.143
$ws.Range("K63").Value = 5043.143
$ws.Range("M63").Value = -4357.143

$ws.Range("H66").Value = 6287.75
$ws.Range("I66").Value = 5043.143
$ws.Range("K66").Value = 25215.715
$ws.Range("M66").Value = -21783.715

$ws.Range("H88").Value = 2732.6875
$ws.Range("I88").Value = 3892.625
$ws.Range("J88").Value = 1572.75
$ws.Range("K88").Value = 3892.625
$ws.Range("L88").Value = 1572.75
$ws.Range("M88").Value = -3486.625
$ws.Range("N88").Value = -2384.75

$ws.Range("H91").Value = 2732.6875
$ws.Range("I91").Value = 3892.625
$ws.Range("J91").Value = 1572.75
$ws.Range("K91").Value = 3892.625
$ws.Range("L91").Value = 1572.75
$ws.Range("M91").Value = -2488.625
$ws.Range("N91").Value = -4380.75

$ws.Range("H124").Value = 42854.715
$ws.Range("J124").Value = 42854.715
$ws.Range("L124").Value = 42854.715
$ws.Range("N124").Value = -52674.715

$ws.Range("H125").Value = 49411.855
$ws.Range("J125").Value = 49411.855
$ws.Range("L125").Value = 49411.855
$ws.Range("N125").Value = -59251.855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1064104.9
$ws.Range("I86").Value = 2429468.2
$ws.Range("K86").Value = 2429468.2
$ws.Range("M86").Value = -2428345.2

$ws.Range("H89").Value = 1064104.9
$ws.Range("I89").Value = 2429468.2
$ws.Range("K89").Value = 12147341
$ws.Range("M89").Value = -12141725

$ws.Range("H103").Value = 50000
$ws.Range("J103").Value = 50000
$ws.Range("L103").Value = 50000
$ws.Range("N103").Value = -52344

$ws.Range("H105").Value = 41301.32
$ws.Range("I105").Value = 46552.137
$ws.Range("J105").Value = 2795.3333
$ws.Range("K105").Value = 46552.137
$ws.Range("L105").Value = 2795.3333
$ws.Range("M105").Value = -44805.137
$ws.Range("N105").Value = -6289.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 2115.1765
$ws.Range("I19").Value = 336.85715
$ws.Range("J19").Value = 3360
$ws.Range("K19").Value = 336.85715
$ws.Range("L19").Value = 3360
$ws.Range("M19").Value = -166.85715
$ws.Range("N19").Value = -3700

$ws.Range("H22").Value = 471.75
$ws.Range("I22").Value = 471.33334
$ws.Range("K22").Value = 471.33334
$ws.Range("M22").Value = -121.33334

$ws.Range("H24").Value = 2115.1765
$ws.Range("I24").Value = 336.85715
$ws.Range("J24").Value = 3360
$ws.Range("K24").Value = 336.85715
$ws.Range("L24").Value = 3360
$ws.Range("M24").Value = -166.85715
$ws.Range("N24").Value = -3700

$ws.Range("H45").Value = 400
$ws.Range("I45").Value = 400
$ws.Range("K45").Value = 400
$ws.Range("M45").Value = 193

$ws.Range("H99").Value = 8186.476
$ws.Range("I99").Value = 9801.700000000001
$ws.Range("K99").Value = 9801.700000000001
$ws.Range("M99").Value = -8303.700000000001

$ws.Range("H126").Value = 8186.476
$ws.Range("I126").Value = 9801.700000000001
$ws.Range("K126").Value = 29405.1
$ws.Range("M126").Value = -26935.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 182154.27
$ws.Range("I26").Value = 333416.66
$ws.Range("K26").Value = 1000249.98
$ws.Range("M26").Value = -999961.98

$ws.Range("H108").Value = 7282.75
$ws.Range("I108").Value = 6894.5713
$ws.Range("K108").Value = 20683.7139
$ws.Range("M108").Value = -17803.7139

$ws.Range("H140").Value = 2919.4546
$ws.Range("I140").Value = 1399
$ws.Range("J140").Value = 3114.3845
$ws.Range("K140").Value = 4197
$ws.Range("L140").Value = 9343.1535
$ws.Range("M140").Value = 983
$ws.Range("N140").Value = -19703.1535

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2380.5
$ws.Range("I102").Value = 1727.3334
$ws.Range("J102").Value = 4340
$ws.Range("K102").Value = 1727.3334
$ws.Range("L102").Value = 4340
$ws.Range("M102").Value = -105.3334
$ws.Range("N102").Value = -7584

$ws.Range("H123").Value = 43999
$ws.Range("J123").Value = 43999
$ws.Range("L123").Value = 43999
$ws.Range("N123").Value = -48899

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 19000
$ws.Range("J5").Value = 19000
$ws.Range("L5").Value = 19000
$ws.Range("N5").Value = -19226

$ws.Range("H40").Value = 23697.348
$ws.Range("I40").Value = 34086.547
$ws.Range("J40").Value = 5804.8335
$ws.Range("K40").Value = 34086.547
$ws.Range("L40").Value = 5804.8335
$ws.Range("M40").Value = -33950.547
$ws.Range("N40").Value = -6076.8335

$ws.Range("H100").Value = 1573.5714
$ws.Range("I100").Value = 1435.8334
$ws.Range("J100").Value = 2400
$ws.Range("K100").Value = 1435.8334
$ws.Range("L100").Value = 2400
$ws.Range("M100").Value = -894.8334
$ws.Range("N100").Value = -3482

$ws.Range("H132").Value = 3865.8857
$ws.Range("I132").Value = 1360.9445
$ws.Range("K132").Value = 4082.8335
$ws.Range("M132").Value = -1552.8335

$ws.Range("H136").Value = 3847.484
$ws.Range("I136").Value = 2899.875
$ws.Range("K136").Value = 8699.625
$ws.Range("M136").Value = -6149.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 92775.42999999999
$ws.Range("J46").Value = 92775.42999999999
$ws.Range("L46").Value = 92775.42999999999
$ws.Range("N46").Value = -93237.42999999999

$ws.Range("H81").Value = 18025.143
$ws.Range("I81").Value = 1552.7142
$ws.Range("J81").Value = 34497.57
$ws.Range("K81").Value = 3105.4284
$ws.Range("L81").Value = 68995.14
$ws.Range("M81").Value = -2044.4284
$ws.Range("N81").Value = -71117.14

$ws.Range("H84").Value = 18025.143
$ws.Range("I84").Value = 1552.7142
$ws.Range("J84").Value = 34497.57
$ws.Range("K84").Value = 15527.142
$ws.Range("L84").Value = 344975.7
$ws.Range("M84").Value = -10223.142
$ws.Range("N84").Value = -355583.7

$ws.Range("H107").Value = 334.2857
$ws.Range("I107").Value = 289.93332
$ws.Range("K107").Value = 869.7999599999999
$ws.Range("M107").Value = 1050.20004

$ws.Range("H113").Value = 1402.7307
$ws.Range("I113").Value = 1031.6111
$ws.Range("K113").Value = 3094.8333
$ws.Range("M113").Value = -924.8333000000002

$ws.Range("H122").Value = 30306132
$ws.Range("J122").Value = 6550
$ws.Range("L122").Value = 19650
$ws.Range("N122").Value = -24550

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H134").Value = 92775.42999999999
$ws.Range("J134").Value = 92775.42999999999
$ws.Range("L134").Value = 278326.29
$ws.Range("N134").Value = -283396.29

$ws.Range("H136").Value = 65286.633
$ws.Range("I136").Value = 14196.781
$ws.Range("K136").Value = 42590.343
$ws.Range("M136").Value = -40040.343
